$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Updated "last refreshed" timestamp
$ws.Range("A1").Value = "Datos actualizados a 11 de Julio de 2020 a las 10:03"

# --- Simple numeric refreshes (country identity / row position unchanged) ---

# Row 6: India
$ws.Range("B6").Value = 822674
$ws.Range("C6").Value = 71
$ws.Range("D6").Value = 516308
$ws.Range("E6").Value = 284214
$ws.Range("G6").Value = 8
$ws.Range("H6").Value = 22152

# Row 7: Rusia
$ws.Range("B7").Value = 720547
$ws.Range("C7").Value = 6611
$ws.Range("D7").Value = 497446
$ws.Range("E7").Value = 211896
$ws.Range("G7").Value = 188
$ws.Range("H7").Value = 11205

# Row 33: Belgica
$ws.Range("B33").Value = 62469
$ws.Range("C33").Value = 112
$ws.Range("D33").Value = 17196
$ws.Range("E33").Value = 35491
$ws.Range("G33").Value = 1
$ws.Range("H33").Value = 9782

# Row 39: Ucrania
$ws.Range("B39").Value = 52843
$ws.Range("C39").Value = 800
$ws.Range("D39").Value = 25661
$ws.Range("E39").Value = 25810
$ws.Range("G39").Value = 27
$ws.Range("H39").Value = 1372

# Row 60: Moldavia
$ws.Range("D60").Value = 12456
$ws.Range("E60").Value = 5833

# Row 68: Chequia
$ws.Range("B68").Value = 13062
$ws.Range("C68").Value = 61
$ws.Range("D68").Value = 8209
$ws.Range("E68").Value = 4501

# Row 117: Estonia
$ws.Range("B117").Value = 2014
$ws.Range("C117").Value = 1
$ws.Range("D117").Value = 1895

# Row 135: Letonia
$ws.Range("B135").Value = 1173
$ws.Range("C135").Value = 8
$ws.Range("E135").Value = 124

# --- Countries whose updated figures move them ahead of a neighbour ---
# Singapur jumps above Portugal: row 41 becomes Singapur (new numbers),
# row 42 becomes Portugal (the old, unchanged Singapur/Portugal row data).
$ws.Range("A41").Value = "Singapur"
$ws.Range("B41").Value = 45783
$ws.Range("C41").Value = 170
$ws.Range("D41").Value = 41780
$ws.Range("E41").Value = 3977
$ws.Range("H41").Value = 26

$ws.Range("A42").Value = "Portugal"
$ws.Range("B42").Value = 45679
$ws.Range("D42").Value = 30350
$ws.Range("E42").Value = 13683
$ws.Range("H42").Value = 1646

# Armenia jumps above Rumania and Nigeria: row 51 becomes Armenia (new
# numbers), row 52 becomes Rumania, row 53 becomes Nigeria (each shifted
# down one row, values unchanged).
$ws.Range("A51").Value = "Armenia"
$ws.Range("B51").Value = 31392
$ws.Range("C51").Value = 489
$ws.Range("D51").Value = 19419
$ws.Range("E51").Value = 11414
$ws.Range("G51").Value = 13
$ws.Range("H51").Value = 559

$ws.Range("A52").Value = "Rumania"
$ws.Range("B52").Value = 31381
$ws.Range("D52").Value = 21129
$ws.Range("E52").Value = 8405
$ws.Range("H52").Value = 1847

$ws.Range("A53").Value = "Nigeria"
$ws.Range("B53").Value = 31323
$ws.Range("D53").Value = 12795
$ws.Range("E53").Value = 17819
$ws.Range("H53").Value = 709

# Eslovaquia jumps above Islandia: row 119 becomes Eslovaquia (new
# numbers), row 120 becomes Islandia (old, unchanged data).
$ws.Range("A119").Value = "Eslovaquia"
$ws.Range("B119").Value = 1893
$ws.Range("C119").Value = 23
$ws.Range("D119").Value = 1493
$ws.Range("E119").Value = 372
$ws.Range("H119").Value = 28

$ws.Range("A120").Value = "Islandia"
$ws.Range("B120").Value = 1886
$ws.Range("D120").Value = 1859
$ws.Range("E120").Value = 17
$ws.Range("H120").Value = 10
